$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The table currently ends at row 72 (A1:C72); we're appending one new
# data row (73) with the latest mod-count snapshot.
$srcRow = 72
$newRow = 73

# 1) Seed row 73's formatting from row 72 (centered alignment, General
#    number format -- the same style already used by every data row) so we
#    don't end up inventing a brand-new cell style just for this row.
$ws.Range("A" + $srcRow + ":C" + $srcRow).Copy()
$ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122)  # xlPasteFormats

# 2) Write the values. Column A holds the date as literal text (exactly
#    like the rest of the column), so build it as a text-formula result
#    and paste back as a value -- that way Excel never reinterprets
#    "2026/01/22" as an actual date serial the way a plain .Value= would.
$ws.Cells.Item($newRow, 1).Formula = "=""2026/01/22"""
$ws.Range("A" + $newRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4163)  # xlPasteValues

$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1155

$excel.CutCopyMode = $false
